$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New consolidated account-statement data: grouped by worker (descending
# period order 1905..1810) instead of grouped by period. Marcos Ricardo's
# "Salario Basico" (col G) is updated to 2000000 for every one of his rows.
$data = @(
    @("CC","73476371","MARCOS RICARDO GARCES VILLA","1905",26666,2000000),
    @("CC","73476371","MARCOS RICARDO GARCES VILLA","1904",53615,2000000),
    @("CC","73476371","MARCOS RICARDO GARCES VILLA","1903",53615,2000000),
    @("CC","73476371","MARCOS RICARDO GARCES VILLA","1902",53615,2000000),
    @("CC","73476371","MARCOS RICARDO GARCES VILLA","1901",53615,2000000),
    @("CC","73476371","MARCOS RICARDO GARCES VILLA","1812",53615,2000000),
    @("CC","73476371","MARCOS RICARDO GARCES VILLA","1810",53615,2000000),
    @("CC","8834880","ALCIDES JESUS CARDENAS LOPEZ","1905",26666,1411788),
    @("CC","8834880","ALCIDES JESUS CARDENAS LOPEZ","1904",40000,1411788),
    @("CC","8834880","ALCIDES JESUS CARDENAS LOPEZ","1903",56472,1411788),
    @("CC","8834880","ALCIDES JESUS CARDENAS LOPEZ","1902",56472,1411788),
    @("CC","8834880","ALCIDES JESUS CARDENAS LOPEZ","1901",56472,1411788),
    @("CC","8834880","ALCIDES JESUS CARDENAS LOPEZ","1812",56472,1411788),
    @("CC","1070822062","DEYMER RAMOS LOPEZ","1905",20833,781242),
    @("CC","1070822062","DEYMER RAMOS LOPEZ","1904",31249,781242),
    @("CC","1070822062","DEYMER RAMOS LOPEZ","1903",31249,781242),
    @("CC","1070822062","DEYMER RAMOS LOPEZ","1902",31249,781242),
    @("CC","1070822062","DEYMER RAMOS LOPEZ","1901",31249,781242),
    @("CC","1070822062","DEYMER RAMOS LOPEZ","1812",31249,781242)
)

$row = 16
foreach ($rec in $data) {
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
    $row = $row + 1
}
